$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.866.45"
$ws.Range("E2").Value = "  +4.71%  "
$ws.Range("D3").Value = "2.273.10"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.11"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.75%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.421"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.86"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0931"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.78%  "
$ws.Range("E12").Value = "  +0.46%  "
$ws.Range("D13").Value = "2.611.81"
$ws.Range("E13").Value = "  +1.96%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.83%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.79"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.807"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "2.277.17"
$ws.Range("E18").Value = "  +1.81%  "
$ws.Range("D19").Value = "43.753.48"
$ws.Range("E19").Value = "  +4.70%  "
$ws.Range("D20").Value = "0.0₃0932"
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("E27").Value = "  +1.82%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "170.82"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  -1.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("E31").Value = "  +2.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0657"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.97%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.38"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("E39").Value = "  -2.42%  "
$ws.Range("E40").Value = "  +4.17%  "
$ws.Range("E41").Value = "  +0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000227"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.67%  "
$ws.Range("E44").Value = "  +1.43%  "
$ws.Range("E45").Value = "  -7.18%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "98.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.60%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("D48").Value = "1.478.47"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("E51").Value = "  -3.25%  "
